# testData_BC.xlsx edit:
#  - devices sheet gets a new data row (row 3) with a device entry
#  - column L widened a bit
#  - devices sheet becomes the active/selected tab (was "signIn")
#  - selection on devices sheet ends on L3
#  - a basic portrait page setup is applied to the devices sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# New row of device data. Write "Available" before "Galaxy S6" so the
# shared-string table grows in the same order as the authored workbook.
$ws.Range("I3").Value = "Available"
$ws.Range("E3").Value = "Galaxy S6"
$ws.Range("K3").Value = "Appium"
$ws.Range("L3").Value = "com.bloomfire.android.perfecto"

# Widen column L (12) from ~29.4 to 33 characters.
$ws.Columns.Item(12).ColumnWidth = 32.14

# Give the sheet a simple portrait page setup.
$ws.PageSetup.Orientation = 1

# Move the selection/active-tab to the devices sheet, ending on L3 -
# this also clears tabSelected on the previously-active "signIn" sheet.
$ws.Range("L3").Select()
